$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark that sits right after the logo image.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Insert "-Capoeira" right after "(ABADÁ" (before "), é uma ...").
$rng = $d.Content
$rng.Find.Execute("(ABADÁ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("-Capoeira")

# 3. "organização sem fins lucrativos" -> "instituição"
$rng2 = $d.Content
$rng2.Find.Execute("organização sem fins lucrativos", $true, $false, $false, $false, $false, $true, 1, $false, "instituição", 2) | Out-Null

# 4. " por " -> " pelo Dr. h. c. " and re-insert the "_GoBack" bookmark right after "pelo".
$rng3 = $d.Content
$rng3.Find.Execute(" por ", $true, $false, $false, $false, $false, $true, 1, $false, " pelo Dr. h. c. ", 2) | Out-Null

$rng4 = $d.Content
$rng4.Find.Execute("criada em 1988 pelo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($rng4.End, $rng4.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
